$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- locate the existing background group ("Group 8") ---
$g = $s.Shapes.Item(1)

# Ungroup it so we can manipulate (and duplicate) the individual leaf-pattern
# pictures that live inside it -- grouped children can't be Copy/Duplicate'd
# directly in this environment.
$ungrouped = $g.Ungroup()

# After ungrouping, the three former group children are now top-level shapes,
# in their original order: Picture 4 (grey backdrop), Picture 6, Picture 7.
$picGrey = $s.Shapes.Item(1)
$pic6 = $s.Shapes.Item(2)
$pic7 = $s.Shapes.Item(3)

# --- create the extra leaf-pattern copies we need ---
# 1) a brand-new picture that will live OUTSIDE the group, to the right,
#    using the original (wider/taller) framing that "Picture 7" used to have.
$picOutside = $pic7.Duplicate()
$picOutside.Name = "Picture 7"
$picOutside.Left = 840.0003937
$picOutside.Top = 283.4598425
$picOutside.Width = 33.72944882
$picOutside.Height = 240.9730709
$picOutside.PictureFormat.CropLeft = 142.09875
$picOutside.PictureFormat.CropTop = -0.76875
$picOutside.PictureFormat.CropRight = 199.17375
$picOutside.PictureFormat.CropBottom = 93.60375

# 2) two further copies of the leaf strip that will stay inside the group.
$pic3New = $pic7.Duplicate()
$pic3New.Name = "Picture 3"
$pic1New = $pic7.Duplicate()
$pic1New.Name = "Picture 1"

# --- resize / reposition the shapes that remain in the group ---

# Picture 6: narrower crop strip, first of four vertical tiles.
$pic6.Left = 726.1377953
$pic6.Top = 3.9147244
$pic6.Width = 17.2136220
$pic6.Height = 144
$pic6.PictureFormat.CropLeft = 142.09875
$pic6.PictureFormat.CropTop = -0.76875
$pic6.PictureFormat.CropRight = 199.17375
$pic6.PictureFormat.CropBottom = 93.60375

# Picture 7 (remaining original) becomes "Picture 2": second vertical tile.
$pic7.Name = "Picture 2"
$pic7.Left = 726.1377953
$pic7.Top = 148.9737008
$pic7.Width = 17.2136220
$pic7.Height = 144
$pic7.PictureFormat.CropLeft = 142.09875
$pic7.PictureFormat.CropTop = -0.76875
$pic7.PictureFormat.CropRight = 199.17375
$pic7.PictureFormat.CropBottom = 93.60375

# Picture 3 (new): third vertical tile.
$pic3New.Left = 726.1377953
$pic3New.Top = 294.0326772
$pic3New.Width = 17.2136220
$pic3New.Height = 144
$pic3New.PictureFormat.CropLeft = 142.09875
$pic3New.PictureFormat.CropTop = -0.76875
$pic3New.PictureFormat.CropRight = 199.17375
$pic3New.PictureFormat.CropBottom = 93.60375

# Picture 1 (new): fourth (shorter) vertical tile.
$pic1New.Left = 726.0714173
$pic1New.Top = 439.0916535
$pic1New.Width = 17.28
$pic1New.Height = 98.9221260
$pic1New.PictureFormat.CropLeft = 142.09875
$pic1New.PictureFormat.CropTop = 91.54875
$pic1New.PictureFormat.CropRight = 199.17375
$pic1New.PictureFormat.CropBottom = 90.36

# --- regroup the five pictures that belong back in the background group ---
$idxGrey = $picGrey.ZOrderPosition
$idxPic6 = $pic6.ZOrderPosition
$idxPic7 = $pic7.ZOrderPosition
$idxPic3 = $pic3New.ZOrderPosition
$idxPic1 = $pic1New.ZOrderPosition
$range = $s.Shapes.Range(@($idxGrey, $idxPic6, $idxPic7, $idxPic3, $idxPic1))
$newGroup = $range.Group()
$newGroup.Name = "Group 8"
